$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.02328126719340038
$ws.Range("J2").Value = 0.02328126719340038
$ws.Range("M2").Value = 10.25883033333333
$ws.Range("N2").Value = 30.776491
$ws.Range("O2").Value = 0.34684992242997
$ws.Range("P2").Value = 0.34684992242997
$ws.Range("Q2").Value = 0.216327955239
$ws.Range("R2").Value = 1.946951597151
$ws.Range("S2").Value = 0.008075105720102328
$ws.Range("T2").Value = 0.008075105720102328

# Row 3
$ws.Range("I3").Value = 0.02328126719340038
$ws.Range("J3").Value = 0.02328126719340038
$ws.Range("O3").Value = 0.1682819529322607
$ws.Range("P3").Value = 0.1682819529322608
$ws.Range("S3").Value = 0.003917817110043189
$ws.Range("T3").Value = 0.003917817110043189

# Row 4
$ws.Range("I4").Value = 0.02328126719340038
$ws.Range("J4").Value = 0.02328126719340038
$ws.Range("M4").Value = 4.130648333333333
$ws.Range("N4").Value = 12.391945
$ws.Range("O4").Value = 0.1396567647041521
$ws.Range("P4").Value = 0.1396567647041521
$ws.Range("Q4").Value = 0.08710298140499999
$ws.Range("R4").Value = 0.783926832645
$ws.Range("S4").Value = 0.003251386454443212
$ws.Range("T4").Value = 0.003251386454443212

# Row 5
$ws.Range("I5").Value = 0.02328126719340038
$ws.Range("J5").Value = 0.02328126719340038
$ws.Range("M5").Value = 2.760918333333333
$ws.Range("N5").Value = 8.282755
$ws.Range("O5").Value = 0.09334634443076846
$ws.Range("P5").Value = 0.09334634443076847
$ws.Range("Q5").Value = 0.05821948489499999
$ws.Range("R5").Value = 0.523975364055
$ws.Range("S5").Value = 0.002173221186219902
$ws.Range("T5").Value = 0.002173221186219903

# Row 6
$ws.Range("I6").Value = 0.02328126719340038
$ws.Range("J6").Value = 0.02328126719340038
$ws.Range("M6").Value = 2.602884
$ws.Range("N6").Value = 7.808651999999999
$ws.Range("O6").Value = 0.08800322104565558
$ws.Range("P6").Value = 0.0880032210456556
$ws.Range("Q6").Value = 0.05488701490799999
$ws.Range("R6").Value = 0.493983134172
$ws.Range("S6").Value = 0.002048826503043783
$ws.Range("T6").Value = 0.002048826503043784

# Row 7
$ws.Range("I7").Value = 0.02328126719340038
$ws.Range("J7").Value = 0.02328126719340038
$ws.Range("M7").Value = 4.846564
$ws.Range("N7").Value = 14.539692
$ws.Range("O7").Value = 0.1638617944571931
$ws.Range("P7").Value = 0.1638617944571932
$ws.Range("Q7").Value = 0.102199495068
$ws.Range("R7").Value = 0.919795455612
$ws.Range("S7").Value = 0.003814910219547967
$ws.Range("T7").Value = 0.003814910219547968

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8846626666666667
$ws.Range("H8").Value = 2.653988
$ws.Range("I8").Value = 0.9767187328065996
$ws.Range("J8").Value = 0.9767187328065997
$ws.Range("M8").Value = 10.25883033333333
$ws.Range("N8").Value = 30.776491
$ws.Range("O8").Value = 0.34684992242997
$ws.Range("P8").Value = 0.34684992242997
$ws.Range("Q8").Value = 9.075604199567557
$ws.Range("R8").Value = 81.680437796108
$ws.Range("S8").Value = 0.3387748167098676
$ws.Range("T8").Value = 0.3387748167098677

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8846626666666667
$ws.Range("H9").Value = 2.653988
$ws.Range("I9").Value = 0.9767187328065996
$ws.Range("J9").Value = 0.9767187328065997
$ws.Range("O9").Value = 0.1682819529322607
$ws.Range("P9").Value = 0.1682819529322608
$ws.Range("Q9").Value = 4.403231195912445
$ws.Range("R9").Value = 39.629080763212
$ws.Range("S9").Value = 0.1643641358222175
$ws.Range("T9").Value = 0.1643641358222176

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8846626666666667
$ws.Range("H10").Value = 2.653988
$ws.Range("I10").Value = 0.9767187328065996
$ws.Range("J10").Value = 0.9767187328065997
$ws.Range("M10").Value = 4.130648333333333
$ws.Range("N10").Value = 12.391945
$ws.Range("O10").Value = 0.1396567647041521
$ws.Range("P10").Value = 0.1396567647041521
$ws.Range("Q10").Value = 3.654230369628889
$ws.Range("R10").Value = 32.88807332666
$ws.Range("S10").Value = 0.1364053782497088
$ws.Range("T10").Value = 0.1364053782497089

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.8846626666666667
$ws.Range("H11").Value = 2.653988
$ws.Range("I11").Value = 0.9767187328065996
$ws.Range("J11").Value = 0.9767187328065997
$ws.Range("M11").Value = 2.760918333333333
$ws.Range("N11").Value = 8.282755
$ws.Range("O11").Value = 0.09334634443076846
$ws.Range("P11").Value = 0.09334634443076847
$ws.Range("Q11").Value = 2.442481375215555
$ws.Range("R11").Value = 21.98233237694
$ws.Range("S11").Value = 0.09117312324454856
$ws.Range("T11").Value = 0.09117312324454858

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.8846626666666667
$ws.Range("H12").Value = 2.653988
$ws.Range("I12").Value = 0.9767187328065996
$ws.Range("J12").Value = 0.9767187328065997
$ws.Range("M12").Value = 2.602884
$ws.Range("N12").Value = 7.808651999999999
$ws.Range("O12").Value = 0.08800322104565558
$ws.Range("P12").Value = 0.0880032210456556
$ws.Range("Q12").Value = 2.302674300464
$ws.Range("R12").Value = 20.724068704176
$ws.Range("S12").Value = 0.08595439454261179
$ws.Range("T12").Value = 0.08595439454261182

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.8846626666666667
$ws.Range("H13").Value = 2.653988
$ws.Range("I13").Value = 0.9767187328065996
$ws.Range("J13").Value = 0.9767187328065997
$ws.Range("M13").Value = 4.846564
$ws.Range("N13").Value = 14.539692
$ws.Range("O13").Value = 0.1638617944571931
$ws.Range("P13").Value = 0.1638617944571932
$ws.Range("Q13").Value = 4.287574232410667
$ws.Range("R13").Value = 38.58816809169601
$ws.Range("S13").Value = 0.1600468842376452
$ws.Range("T13").Value = 0.1600468842376452

Write-Output "Updated TPM-derived values for Wnt7b-Lrp5 LR pairs (rows 2-13)"
